$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-converted to a number by Excel (losing formatting like trailing zeros).
$textRows = @(5,6,7,8,9,10,11,15,16,18,19,20,22,23,24,25,26,29,30,32,33,34,35,37,38,39,41,42,44,45,46,48,49,50)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range('D2').Value = '60.328.27'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '2.629.85'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '522.96'
$ws.Range('E5').Value = '  +1.78%  '
$ws.Range('D6').Value = '151.07'
$ws.Range('E6').Value = '  -1.45%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.573'
$ws.Range('E8').Value = '  -4.09%  '
$ws.Range('D9').Value = '6.40'
$ws.Range('E9').Value = '  -3.83%  '
$ws.Range('D10').Value = '0.106'
$ws.Range('E10').Value = '  +2.38%  '
$ws.Range('D11').Value = '0.344'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '3.091.20'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').Value = '60.366.58'
$ws.Range('E14').Value = '  -0.09%  '
$ws.Range('D15').Value = '21.50'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').Value = '0.0000139'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '2.634.57'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').Value = '4.66'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').Value = '347.20'
$ws.Range('E19').Value = '  -3.12%  '
$ws.Range('D20').Value = '10.47'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').Value = '0.995'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').Value = '60.93'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '0.421'
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('D25').Value = '0.164'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = '0.0₃0830'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '6.05'
$ws.Range('E30').Value = '  +2.08%  '
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').Value = '19.05'
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').Value = '150.26'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '4.00'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '0.900'
$ws.Range('E35').Value = '  -0.69%  '
$ws.Range('E36').Value = '  -2.19%  '
$ws.Range('D37').Value = '0.879'
$ws.Range('E37').Value = '  +4.38%  '
$ws.Range('D38').Value = '36.67'
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('D39').Value = '1.45'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('E40').Value = '  -1.71%  '
$ws.Range('D41').Value = '290.27'
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('D42').Value = '0.634'
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').Value = '19.82'
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('D46').Value = '0.0552'
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').Value = '4.78'
$ws.Range('E48').Value = '  -3.63%  '
$ws.Range('D49').Value = '10.40'
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('D50').Value = '18.84'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('D51').Value = '1.968.71'
$ws.Range('E51').Value = '  -1.07%  '
